$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Convert row 2 numeric-looking text cells (C2:G2) to real numbers
$ws.Range("C2").Value = 14.8881
$ws.Range("D2").Value = 120.7855
$ws.Range("E2").Value = 4663
$ws.Range("F2").Value = 92
$ws.Range("G2").Value = 1000

# Add new row 3 data
$ws.Range("A3").Value = $true
$ws.Range("B3").Value = "bry"
$ws.Range("C3").Value = 14.9
$ws.Range("D3").Value = 120.78
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 3
$ws.Range("G3").Value = 2
$ws.Range("H3").Value = "d"
